$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Recommandations")
$ws1.Cells.Item(2, 4).Value = 83370.3
$ws1.Cells.Item(2, 5).Value = 108.68
$ws1.Cells.Item(3, 4).Value = 78215
$ws1.Cells.Item(4, 3).Value = 76
$ws1.Cells.Item(4, 4).Value = 71225
$ws1.Cells.Item(5, 4).Value = 65283.98
$ws1.Cells.Item(5, 5).Value = 653.39
$ws1.Cells.Item(6, 4).Value = 61355
$ws1.Cells.Item(6, 5).Value = 745
$ws1.Cells.Item(7, 4).Value = 60130
$ws1.Cells.Item(7, 5).Value = 540
$ws1.Cells.Item(8, 4).Value = 59095
$ws1.Cells.Item(8, 5).Value = 655
$ws1.Cells.Item(9, 4).Value = 54720
$ws1.Cells.Item(9, 5).Value = 545
$ws1.Cells.Item(10, 4).Value = 49025
$ws1.Cells.Item(11, 4).Value = 42790
$ws1.Cells.Item(11, 5).Value = 510
$ws1.Cells.Item(12, 4).Value = 36645.03
$ws1.Cells.Item(12, 5).Value = 357.08
$ws1.Cells.Item(13, 4).Value = 32648.83
$ws1.Cells.Item(13, 5).Value = 366.76
$ws1.Cells.Item(14, 4).Value = 23086.65
$ws1.Cells.Item(14, 5).Value = 284.64
$ws1.Cells.Item(16, 4).Value = 15139.44
$ws1.Cells.Item(16, 5).Value = 169.24
$ws1.Cells.Item(17, 4).Value = 13884.78
$ws1.Cells.Item(17, 5).Value = 172.36
$ws1.Cells.Item(18, 4).Value = 11955.81
$ws1.Cells.Item(18, 5).Value = 149.2
$ws1.Cells.Item(19, 4).Value = 11890.56
$ws1.Cells.Item(19, 5).Value = 126.83
$ws1.Cells.Item(20, 4).Value = 11394.81
$ws1.Cells.Item(20, 5).Value = 137.16
$ws1.Cells.Item(21, 4).Value = 11200.04
$ws1.Cells.Item(21, 5).Value = 104.03
$ws1.Cells.Item(22, 4).Value = 10965.62
$ws1.Cells.Item(22, 5).Value = 120.24
$ws1.Cells.Item(23, 4).Value = 10776.74
$ws1.Cells.Item(23, 5).Value = 118.17
$ws1.Cells.Item(24, 4).Value = 9675.31
$ws1.Cells.Item(24, 5).Value = 104.15
$ws1.Cells.Item(25, 4).Value = 9460.05
$ws1.Cells.Item(25, 5).Value = 96.69
$ws1.Cells.Item(26, 3).Value = 7
$ws1.Cells.Item(26, 4).Value = 4200
$ws1.Cells.Item(26, 5).Value = 475
$ws1.Cells.Item(28, 2).Value = 20
$ws1.Cells.Item(28, 4).Value = 97.88
$ws1.Cells.Item(28, 5).Value = 7.44
$ws1.Cells.Item(28, 7).Value = "➖ Neutre"
$ws1.Cells.Item(30, 2).Value = 19
$ws1.Cells.Item(30, 4).Value = 63.78
$ws1.Cells.Item(31, 1).Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Cells.Item(31, 2).Value = 24
$ws1.Cells.Item(31, 3).Value = 20
$ws1.Cells.Item(31, 4).Value = 47.27
$ws1.Cells.Item(31, 5).Value = -3.7
$ws1.Cells.Item(31, 7).Value = "✅ Renforcer"
$ws1.Cells.Item(32, 1).Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws1.Cells.Item(32, 2).Value = 12
$ws1.Cells.Item(32, 3).Value = 3
$ws1.Cells.Item(32, 4).Value = 45.14
$ws1.Cells.Item(32, 5).Value = -3.36
$ws1.Cells.Item(32, 7).Value = "➖ Neutre"
$ws1.Cells.Item(33, 1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Cells.Item(33, 2).Value = 25
$ws1.Cells.Item(33, 3).Value = 20
$ws1.Cells.Item(33, 4).Value = 42.84
$ws1.Cells.Item(33, 5).Value = 6.25
$ws1.Cells.Item(33, 7).Value = "✅ Renforcer"
$ws1.Cells.Item(34, 1).Value = "BICI CI (BICC)"
$ws1.Cells.Item(34, 2).Value = 8
$ws1.Cells.Item(34, 3).Value = 3
$ws1.Cells.Item(34, 4).Value = 42.51
$ws1.Cells.Item(34, 5).Value = 7.48
$ws1.Cells.Item(34, 7).Value = "✅ Renforcer"
$ws1.Cells.Item(35, 1).Value = "BANK OF AFRICA SENEGAL (BOAS)"
$ws1.Cells.Item(35, 2).Value = 10
$ws1.Cells.Item(35, 3).Value = 5
$ws1.Cells.Item(35, 4).Value = 39.73
$ws1.Cells.Item(35, 5).Value = 3.92
$ws1.Cells.Item(35, 7).Value = "Non évalué"
$ws1.Cells.Item(36, 2).Value = 13
$ws1.Cells.Item(36, 4).Value = 36.42
$ws1.Cells.Item(37, 1).Value = "UNIWAX CI (UNXC)"
$ws1.Cells.Item(37, 2).Value = 21
$ws1.Cells.Item(37, 3).Value = 16
$ws1.Cells.Item(37, 4).Value = 36.08
$ws1.Cells.Item(37, 5).Value = -6.86
$ws1.Cells.Item(37, 7).Value = "⚠️ Risque de décrochage"
$ws1.Cells.Item(38, 1).Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Cells.Item(38, 2).Value = 9
$ws1.Cells.Item(38, 3).Value = 3
$ws1.Cells.Item(38, 4).Value = 32.56
$ws1.Cells.Item(38, 5).Value = 5
$ws1.Cells.Item(39, 1).Value = "TOTALENERGIES MARKETING SN (TTLS)"
$ws1.Cells.Item(39, 2).Value = 12
$ws1.Cells.Item(39, 3).Value = 7
$ws1.Cells.Item(39, 4).Value = 31.07
$ws1.Cells.Item(39, 5).Value = -3.41
$ws1.Cells.Item(39, 7).Value = "Non évalué"
$ws1.Cells.Item(40, 7).Value = "➖ Neutre"
$ws1.Cells.Item(41, 1).Value = "CFAO MOTORS CI (CFAC)"
$ws1.Cells.Item(41, 2).Value = 11
$ws1.Cells.Item(41, 3).Value = 7
$ws1.Cells.Item(41, 4).Value = 24.97
$ws1.Cells.Item(41, 5).Value = 6.98
$ws1.Cells.Item(42, 1).Value = "SMB CI (SMBC)"
$ws1.Cells.Item(42, 2).Value = 9
$ws1.Cells.Item(42, 4).Value = 23.09
$ws1.Cells.Item(42, 5).Value = -2.07
$ws1.Cells.Item(42, 7).Value = "Non évalué"
$ws1.Cells.Item(43, 1).Value = "SOLIBRA CI (SLBC)"
$ws1.Cells.Item(43, 2).Value = 14
$ws1.Cells.Item(43, 3).Value = 11
$ws1.Cells.Item(43, 4).Value = 22.81
$ws1.Cells.Item(43, 5).Value = -7.35
$ws1.Cells.Item(43, 7).Value = "👀 À surveiller"
$ws1.Cells.Item(44, 1).Value = "BANK OF AFRICA ML (BOAM)"
$ws1.Cells.Item(44, 2).Value = 13
$ws1.Cells.Item(44, 3).Value = 9
$ws1.Cells.Item(44, 4).Value = 18.35
$ws1.Cells.Item(44, 5).Value = 3.2
$ws1.Cells.Item(44, 7).Value = "Non évalué"
$ws1.Cells.Item(45, 1).Value = "BANK OF AFRICA BN (BOAB)"
$ws1.Cells.Item(45, 2).Value = 4
$ws1.Cells.Item(45, 3).Value = 3
$ws1.Cells.Item(45, 4).Value = 14.2
$ws1.Cells.Item(45, 5).Value = -1.9
$ws1.Cells.Item(45, 7).Value = "➖ Neutre"
$ws1.Cells.Item(46, 1).Value = "BERNABE CI (BNBC)"
$ws1.Cells.Item(46, 2).Value = 22
$ws1.Cells.Item(46, 3).Value = 21
$ws1.Cells.Item(46, 4).Value = 13.54
$ws1.Cells.Item(46, 5).Value = 5.05
$ws1.Cells.Item(46, 7).Value = "⚠️ Risque de décrochage"
$ws1.Cells.Item(49, 1).Value = "SAFCA CI (SAFC)"
$ws1.Cells.Item(49, 2).Value = 11
$ws1.Cells.Item(49, 3).Value = 8
$ws1.Cells.Item(49, 4).Value = 8.25
$ws1.Cells.Item(49, 5).Value = 6.47
$ws1.Cells.Item(49, 7).Value = "👀 À surveiller"
$ws1.Cells.Item(50, 1).Value = "ORAGROUP TOGO (ORGT)"
$ws1.Cells.Item(50, 2).Value = 7
$ws1.Cells.Item(50, 3).Value = 6
$ws1.Cells.Item(50, 4).Value = 8
$ws1.Cells.Item(50, 5).Value = -4.72
$ws1.Cells.Item(50, 7).Value = "➖ Neutre"
$ws1.Cells.Item(51, 1).Value = "SICOR CI (SICC)"
$ws1.Cells.Item(51, 2).Value = 9
$ws1.Cells.Item(51, 3).Value = 9
$ws1.Cells.Item(51, 4).Value = 5.92
$ws1.Cells.Item(51, 5).Value = 6.97
$ws1.Cells.Item(52, 1).Value = "VIVO ENERGY CI (SHEC)"
$ws1.Cells.Item(52, 2).Value = 10
$ws1.Cells.Item(52, 3).Value = 10
$ws1.Cells.Item(52, 4).Value = 5.73
$ws1.Cells.Item(52, 5).Value = -2.12
$ws1.Cells.Item(52, 6).Value = "🟢 Achat"
$ws1.Cells.Item(52, 7).Value = "Non évalué"
$ws1.Cells.Item(53, 1).Value = "ONATEL BF (ONTBF)"
$ws1.Cells.Item(53, 2).Value = 6
$ws1.Cells.Item(53, 3).Value = 11
$ws1.Cells.Item(53, 4).Value = 4.49
$ws1.Cells.Item(53, 5).Value = 3.95
$ws1.Cells.Item(53, 7).Value = "➖ Neutre"
$ws1.Cells.Item(54, 1).Value = "SOGB CI (SOGC)"
$ws1.Cells.Item(54, 2).Value = 10
$ws1.Cells.Item(54, 3).Value = 6
$ws1.Cells.Item(54, 4).Value = 3.94
$ws1.Cells.Item(54, 5).Value = 2.28
$ws1.Cells.Item(55, 1).Value = "ORANGE COTE D'IVOIRE (ORAC)"
$ws1.Cells.Item(55, 2).Value = 17
$ws1.Cells.Item(55, 3).Value = 13
$ws1.Cells.Item(55, 4).Value = 3.59
$ws1.Cells.Item(55, 5).Value = 3.52
$ws1.Cells.Item(55, 7).Value = "Non évalué"
$ws1.Cells.Item(58, 1).Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws1.Cells.Item(58, 2).Value = 1
$ws1.Cells.Item(58, 4).Value = -7.75
$ws1.Cells.Item(58, 5).Value = -1.14
$ws1.Cells.Item(58, 7).Value = "Non évalué"
$ws1.Cells.Item(59, 1).Value = "SONATEL SN (SNTS)"
$ws1.Cells.Item(59, 2).Value = 4
$ws1.Cells.Item(59, 3).Value = 9
$ws1.Cells.Item(59, 4).Value = -8.76
$ws1.Cells.Item(59, 5).Value = -2.89
$ws1.Cells.Item(59, 7).Value = "👀 À surveiller"
$ws1.Cells.Item(60, 1).Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Cells.Item(60, 2).Value = 0
$ws1.Cells.Item(60, 3).Value = 4
$ws1.Cells.Item(60, 4).Value = -8.93
$ws1.Cells.Item(60, 5).Value = -2.63
$ws1.Cells.Item(64, 1).Value = "BANK OF AFRICA BF (BOABF)"
$ws1.Cells.Item(64, 2).Value = 9
$ws1.Cells.Item(64, 3).Value = 13
$ws1.Cells.Item(64, 4).Value = -12.42
$ws1.Cells.Item(64, 5).Value = -2.91
$ws1.Cells.Item(64, 7).Value = "➖ Neutre"
$ws1.Cells.Item(65, 2).Value = 18
$ws1.Cells.Item(65, 4).Value = -15.68
$ws1.Cells.Item(65, 5).Value = 3.7
$ws1.Cells.Item(65, 7).Value = "👀 À surveiller"
$ws1.Cells.Item(70, 1).Value = "CORIS BANK INTERNATIONAL (CBIBF)"
$ws1.Cells.Item(70, 2).Value = 7
$ws1.Cells.Item(70, 3).Value = 16
$ws1.Cells.Item(70, 4).Value = -34.77
$ws1.Cells.Item(70, 5).Value = 3.93
$ws1.Cells.Item(70, 7).Value = "👀 À surveiller"
$ws1.Cells.Item(71, 1).Value = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$ws1.Cells.Item(71, 2).Value = 2
$ws1.Cells.Item(71, 3).Value = 13
$ws1.Cells.Item(71, 4).Value = -36.38
$ws1.Cells.Item(71, 5).Value = -1.28
$ws1.Cells.Item(71, 7).Value = "➖ Neutre"
$ws1.Cells.Item(72, 2).Value = 5
$ws1.Cells.Item(72, 4).Value = -52.3

$ws2 = $wb.Worksheets.Item("Top_YTD")
$ws2.Cells.Item(3, 2).Value = 153.27
$ws2.Cells.Item(5, 2).Value = 78.56
$ws2.Cells.Item(6, 2).Value = 54.39
$ws2.Cells.Item(7, 1).Value = "BICI CI (BICC)"
$ws2.Cells.Item(7, 2).Value = 49.27
$ws2.Cells.Item(8, 1).Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws2.Cells.Item(8, 2).Value = 49.16
$ws2.Cells.Item(9, 1).Value = "BANK OF AFRICA SENEGAL (BOAS)"
$ws2.Cells.Item(9, 2).Value = 46.75
$ws2.Cells.Item(10, 2).Value = 40.36
$ws2.Cells.Item(11, 1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws2.Cells.Item(11, 2).Value = 39.45
